{"js": "// The commit underlines the second \"Dada la cadena 'Donostia'\" bullet\n// (the one immediately followed by \"Mostrar que car\u00e1cter contiene en la\n// posici\u00f3n 2\") and leaves Word's \"last edit location\" bookmark\n// (`_GoBack`) around the text \"Comprobar si contiene la letra e y en que\n// posici\u00f3n\" \u2014 both are applied below using the Word JS API.\n\nconst body = context.document.body;\n\n// --- 1. Underline the 2nd \"Dada la cadena ...\" paragraph -----------------\nconst donostiaResults = body.search(\"Dada la cadena\", { matchCase: false });\ndonostiaResults.load(\"items\");\nawait context.sync();\n\nif (donostiaResults.items.length < 2) {\n  throw new Error(\"Expected to find 2 occurrences of 'Dada la cadena'\");\n}\n\nconst secondDonostia = donostiaResults.items[1];\nconst donostiaParagraph = secondDonostia.paragraphs.getFirst();\nconst donostiaRange = donostiaParagraph.getRange();\ndonostiaRange.font.underline = Word.UnderlineType.single;\n\n// --- 2. Re-create the \"_GoBack\" bookmark around the edited text ----------\nconst bookmarkResults = body.search(\n  \"Comprobar si contiene la letra e y en que posici\u00f3n\",\n  { matchCase: false }\n);\nbookmarkResults.load(\"items\");\nawait context.sync();\n\nif (bookmarkResults.items.length < 1) {\n  throw new Error(\"Could not find the bookmark anchor text\");\n}\n\nbookmarkResults.items[0].insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# The commit underlines the second \"Dada la cadena 'Donostia'\" bullet\n# (the one immediately followed by \"Mostrar que car\u00e1cter contiene en la\n# posici\u00f3n 2\") and leaves Word's \"last edit location\" bookmark\n# (`_GoBack`) around the text \"Comprobar si contiene la letra e y en que\n# posici\u00f3n\" \u2014 both are applied below using the Word COM object model.\n\n$d = $word.ActiveDocument\n\n# --- 1. Underline the 2nd \"Dada la cadena ...\" paragraph -----------------\n$paras = $d.Paragraphs\n$paraCount = $paras.Count\n\n$matchingIndexes = @()\nfor ($i = 1; $i -le $paraCount; $i++) {\n    if ($paras.Item($i).Range.Text -like \"*Dada la cadena*\") {\n        $matchingIndexes += $i\n    }\n}\n\n$targetParaIndex = $matchingIndexes[1]\n$paras.Item($targetParaIndex).Range.Font.Underline = 1\n\n# --- 2. Re-create the \"_GoBack\" bookmark around the edited text ----------\n$bookmarkRange = $d.Content\n$null = $bookmarkRange.Find.Execute(\"Comprobar si contiene la letra e y en que posici\u00f3n\")\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n"}
